$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 into the new I1:J1 header cells
$ws.Range("H1").Copy($ws.Range("I1:J1"))

# Set header text for new columns I (I0) and J (IF) after the style copy so the
# values are not overwritten by the copied content of H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for I2:J22
$data = @(
    @(5, 6),
    @(4, 6),
    @(10, 11),
    @(8, 8),
    @(7, 8),
    @(3, 4),
    @(7, 8),
    @(4, 4),
    @(8, 8),
    @(9, 9),
    @(4, 5),
    @(4, 4),
    @(4, 5),
    @(4, 4),
    @(5, 6),
    @(4, 6),
    @(5, 6),
    @(4, 5),
    @(8, 8),
    @(5, 6),
    @(6, 7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
